# Add an "Interval" column to the Test Schedules Import sheet.
# A new column is inserted before the existing "Lead Weeks (override)"
# column (previously column G), pushing the old G/H/I columns (Lead
# Weeks (override), Timezone (override), Notes) to H/I/J. The new
# column G gets the header "Interval" and small integer values per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G, shifting G:I -> H:J and extending formatting.
$ws.Columns("G:G").Insert()

# Header for the new column.
$ws.Range("G1").Value = "Interval"

# Interval values for the data rows.
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 1

# Match the saved selection from the authored edit.
$ws.Range("G11").Select()
